$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.327.24"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").Value = "3.356.47"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.09"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.19"
$ws.Range("E6").Value = "  +0.53%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.589"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("E9").Value = "  +3.33%  "
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "47.98"
$ws.Range("E11").Value = "  +5.61%  "
$ws.Range("E12").Value = "  +1.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "686.51"
$ws.Range("E13").Value = "  +1.88%  "
$ws.Range("D14").Value = "3.908.02"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").Value = "68.382.21"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.120"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.351.89"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.50"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.21"
$ws.Range("E20").Value = "  +2.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.896"
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.48"
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.96"
$ws.Range("E23").Value = "  -0.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "100.15"
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("E25").Value = "  +1.69%  "
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("E27").Value = "  +2.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.01"
$ws.Range("E28").Value = "  -2.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.53"
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("E30").Value = "  -4.88%  "
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "551.22"
$ws.Range("E32").Value = "  -3.70%  "
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "58.13"
$ws.Range("E34").Value = "  +2.59%  "
$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "3.717.63"
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.36"
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.139"
$ws.Range("E38").Value = "  +5.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.72"
$ws.Range("E39").Value = "  +0.64%  "
$ws.Range("E40").Value = "  +2.17%  "
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("D42").Value = "0.0₃0673"
$ws.Range("E42").Value = "  +0.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.335"
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.25"
$ws.Range("E44").Value = "  -1.88%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0412"
$ws.Range("E45").Value = "  +1.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.64"
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.129"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("E48").Value = "  -0.12%  "
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.66"
$ws.Range("E50").Value = "  +1.99%  "
$ws.Range("E51").Value = "  -1.24%  "
